$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 45181 = 2023-09-12) for every
# data row (2 through 490). Bump it by one day to 45182 (2023-09-13) for all rows.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 490 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45182
